$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.706.09'
$ws.Range("E2").Value = '  -7.14%  '
$ws.Range("D3").Value = '3.307.04'
$ws.Range("E3").Value = '  -8.06%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = "'184.45"
$ws.Range("E5").Value = '  -11.38%  '
$ws.Range("D6").Value = "'520.50"
$ws.Range("E6").Value = '  -8.50%  '
$ws.Range("E7").Value = '  -2.19%  '
$ws.Range("D8").Value = '3.298.00'
$ws.Range("E8").Value = '  -8.12%  '
$ws.Range("D10").Value = "'0.624"
$ws.Range("E10").Value = '  -8.41%  '
$ws.Range("D11").Value = "'60.37"
$ws.Range("E11").Value = '  -5.31%  '
$ws.Range("E12").Value = '  -11.00%  '
$ws.Range("E13").Value = '  -8.94%  '
$ws.Range("E14").Value = '  -9.36%  '
$ws.Range("D15").Value = '3.817.97'
$ws.Range("E15").Value = '  -8.24%  '
$ws.Range("D16").Value = "'0.119"
$ws.Range("E16").Value = '  -4.95%  '
$ws.Range("D17").Value = '3.297.54'
$ws.Range("E17").Value = '  -8.20%  '
$ws.Range("D18").Value = "'17.61"
$ws.Range("E18").Value = '  -8.56%  '
$ws.Range("D19").Value = '63.575.11'
$ws.Range("E19").Value = '  -6.83%  '
$ws.Range("D20").Value = "'11.10"
$ws.Range("E20").Value = '  -9.22%  '
$ws.Range("D21").Value = "'0.955"
$ws.Range("E21").Value = '  -10.48%  '
$ws.Range("D22").Value = "'374.71"
$ws.Range("E22").Value = '  -7.89%  '
$ws.Range("D23").Value = "'11.45"
$ws.Range("E23").Value = '  -7.05%  '
$ws.Range("D24").Value = "'80.72"
$ws.Range("E24").Value = '  -4.87%  '
$ws.Range("D25").Value = "'3.70"
$ws.Range("E25").Value = '  -11.20%  '
$ws.Range("D26").Value = "'3.92"
$ws.Range("E26").Value = '  +1.74%  '
$ws.Range("E27").Value = '  -3.37%  '
$ws.Range("E28").Value = '  -7.34%  '
$ws.Range("D29").Value = "'11.50"
$ws.Range("E29").Value = '  -7.96%  '
$ws.Range("D30").Value = "'8.40"
$ws.Range("E30").Value = '  -7.99%  '
$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").Value = "'652.88"
$ws.Range("E31").Value = '  -11.35%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = "'28.68"
$ws.Range("E32").Value = '  -9.27%  '
$ws.Range("D33").Value = "'6.78"
$ws.Range("E33").Value = '  -10.67%  '
$ws.Range("D34").Value = "'11.25"
$ws.Range("E34").Value = '  -7.37%  '
$ws.Range("D35").Value = "'59.78"
$ws.Range("E35").Value = '  -6.36%  '
$ws.Range("E36").Value = '  -6.26%  '
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("D38").Value = "'0.395"
$ws.Range("E38").Value = '  -6.43%  '
$ws.Range("D39").Value = "'36.58"
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = '  +0.29%  '
$ws.Range("D41").Value = '2.994.19'
$ws.Range("E41").Value = '  -5.21%  '
$ws.Range("E42").Value = '  -4.54%  '
$ws.Range("D43").Value = '0.0₃0659'
$ws.Range("E43").Value = '  -11.70%  '
$ws.Range("D44").Value = "'2.74"
$ws.Range("E44").Value = '  -15.97%  '
$ws.Range("D45").Value = "'2.45"
$ws.Range("E45").Value = '  -5.84%  '
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = "'0.0392"
$ws.Range("E46").Value = '  -5.14%  '
$ws.Range("B47").Value = 'WEMIXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").Value = "'2.61"
$ws.Range("E47").Value = '  -6.00%  '
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").Value = "'2.83"
$ws.Range("E48").Value = '  +4.50%  '
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").Value = "'0.126"
$ws.Range("E49").Value = '  -3.91%  '
$ws.Range("B50").Value = 'ApeXProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D50").Value = "'2.98"
$ws.Range("E50").Value = '  -4.14%  '
$ws.Range("D51").Value = "'2.47"
$ws.Range("E51").Value = '  -22.56%  '
